# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row numbers (in sheet "展览") -> new value for column F
$sheet1Updates = @(
    @{ Row = 2;  New = 119  },
    @{ Row = 3;  New = 255  },
    @{ Row = 4;  New = 132  },
    @{ Row = 5;  New = 1669 },
    @{ Row = 6;  New = 1477 },
    @{ Row = 7;  New = 270  },
    @{ Row = 9;  New = 432  },
    @{ Row = 10; New = 116  }
)

# Row numbers (in sheet "全部类型") -> new value for column F
$sheet4Updates = @(
    @{ Row = 2;  New = 119  },
    @{ Row = 3;  New = 255  },
    @{ Row = 4;  New = 132  },
    @{ Row = 5;  New = 1669 },
    @{ Row = 6;  New = 1477 },
    @{ Row = 7;  New = 270  },
    @{ Row = 10; New = 432  },
    @{ Row = 11; New = 116  }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $sheet1Updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.New
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $sheet4Updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.New
}
